$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Actualización del periodo reportado (cambio de fracciones e histórico):
# de 2do trimestre 2022 (abr-jun) a 3er trimestre 2022 (jul-sep)
$ws.Range("B8").Value = 44743   # Fecha de inicio del periodo que se informa
$ws.Range("C8").Value = 44834   # Fecha de término del periodo que se informa
$ws.Range("E8").Value = 44834   # Fecha de elaboración
$ws.Range("H8").Value = 44844   # Fecha de validación
$ws.Range("I8").Value = 44844   # Fecha de actualización

# Ajuste de la vista de la hoja
$ws.Range("B16").Select() | Out-Null
